# Applies the "Handles float input without breaking stuff" marksheet update.
# - Updates the summary score block (rows 10-12)
# - Normalizes C11 to a numeric -1 instead of text "-1"
# - Fills in "Student Ans" values (col A, and col D for a couple of rows) with
#   correct/incorrect/not-attempted formatting copied from existing styled cells
# - Removes the third (G/H) "Student Ans / Correct Ans" block entirely, and
#   removes most of the second (D/E) block except rows 16-18

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$xlPasteFormats = -4122

function Copy-Style($fromAddr, $toAddr) {
    $ws.Range($fromAddr).Copy() | Out-Null
    $ws.Range($toAddr).PasteSpecial($xlPasteFormats) | Out-Null
}

# ---- Summary block (rows 10-12) ----
# Give A10/A11/A12 the same "mtitleStyle" look already used by A9.
Copy-Style "A9" "A10"
Copy-Style "A9" "A11"
Copy-Style "A9" "A12"

$ws.Range("A10").Value = "No."
$ws.Range("B10").Value = 14
$ws.Range("C10").Value = 1
$ws.Range("D10").Value = 13
$ws.Range("E10").Value = 28

$ws.Range("A11").Value = "Marking"
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -1
$ws.Range("D11").Value = 0

$ws.Range("A12").Value = "Total"
$ws.Range("B12").Value = 56
$ws.Range("C12").Value = -1
$ws.Range("E12").Value = "55/112"

# ---- Drop the third Student Ans / Correct Ans block (columns G:H) ----
$ws.Range("G15:H40").Clear() | Out-Null

# ---- Row 15 header stays as-is for columns A, B, D, E ----

# ---- Row 16: fill in student answer (correct) in column A ----
Copy-Style "B10" "A16"
$ws.Range("A16").Value = "Option A"

# ---- Row 17: student left col A blank (not attempted); fill col D (correct) ----
Copy-Style "B10" "D17"
$ws.Range("D17").Value = "Option C"

# ---- Row 18: student left col A blank (not attempted); fill col D (correct) ----
Copy-Style "B10" "D18"
$ws.Range("D18").Value = "Option D"

# ---- Row 19: fill col A (correct); clear D/E entirely ----
Copy-Style "B10" "A19"
$ws.Range("A19").Value = "Option C"
$ws.Range("D19:E19").Clear() | Out-Null

# ---- Row 20: col A stays blank (not attempted); clear D/E entirely ----
$ws.Range("D20:E20").Clear() | Out-Null

# ---- Row 21: fill col A (correct); clear D/E entirely ----
Copy-Style "B10" "A21"
$ws.Range("A21").Value = "Option C"
$ws.Range("D21:E21").Clear() | Out-Null

# ---- Row 22: fill col A (correct); clear D/E ----
Copy-Style "B10" "A22"
$ws.Range("A22").Value = "Option D"
$ws.Range("D22:E22").Clear() | Out-Null

# ---- Row 23: fill col A (correct); clear D/E ----
Copy-Style "B10" "A23"
$ws.Range("A23").Value = "Option D"
$ws.Range("D23:E23").Clear() | Out-Null

# ---- Row 24: col A stays blank; clear D/E ----
$ws.Range("D24:E24").Clear() | Out-Null

# ---- Row 25: col A stays blank; clear D/E ----
$ws.Range("D25:E25").Clear() | Out-Null

# ---- Row 26: col A stays blank; clear D/E ----
$ws.Range("D26:E26").Clear() | Out-Null

# ---- Row 27: fill col A (correct); clear D/E ----
Copy-Style "B10" "A27"
$ws.Range("A27").Value = "Option A"
$ws.Range("D27:E27").Clear() | Out-Null

# ---- Row 28: fill col A (correct); clear D/E ----
Copy-Style "B10" "A28"
$ws.Range("A28").Value = "Option D"
$ws.Range("D28:E28").Clear() | Out-Null

# ---- Row 29: fill col A (correct); clear D/E ----
Copy-Style "B10" "A29"
$ws.Range("A29").Value = "Option D"
$ws.Range("D29:E29").Clear() | Out-Null

# ---- Row 30: col A stays blank; clear D/E ----
$ws.Range("D30:E30").Clear() | Out-Null

# ---- Row 31: fill col A (correct); clear D/E ----
Copy-Style "B10" "A31"
$ws.Range("A31").Value = "Option D"
$ws.Range("D31:E31").Clear() | Out-Null

# ---- Row 32: fill col A (correct); clear D/E ----
Copy-Style "B10" "A32"
$ws.Range("A32").Value = "Option C"
$ws.Range("D32:E32").Clear() | Out-Null

# ---- Row 33: fill col A (correct); clear D/E ----
Copy-Style "B10" "A33"
$ws.Range("A33").Value = "Option D"
$ws.Range("D33:E33").Clear() | Out-Null

# ---- Row 34: col A stays blank; clear D/E ----
$ws.Range("D34:E34").Clear() | Out-Null

# ---- Row 35: col A stays blank; clear D/E ----
$ws.Range("D35:E35").Clear() | Out-Null

# ---- Row 36: fill col A (INCORRECT - student picked D, correct was A); clear D/E ----
Copy-Style "C10" "A36"
$ws.Range("A36").Value = "Option D"
$ws.Range("D36:E36").Clear() | Out-Null

# ---- Row 37: col A stays blank; clear D/E ----
$ws.Range("D37:E37").Clear() | Out-Null

# ---- Row 38: col A stays blank; clear D/E ----
$ws.Range("D38:E38").Clear() | Out-Null

# ---- Row 39: fill col A (correct); clear D/E ----
Copy-Style "B10" "A39"
$ws.Range("A39").Value = "Option D"
$ws.Range("D39:E39").Clear() | Out-Null

# ---- Row 40: col A stays blank; clear D/E ----
$ws.Range("D40:E40").Clear() | Out-Null
